$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark all test case Runmode values to "Y" (yes)
$ws.Range("C2").Value = "Y"
$ws.Range("C3").Value = "Y"
